$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.618.81"
$ws.Range("E2").Value = "  +2.96%  "

# Row 3
$ws.Range("D3").Value = "3.385.02"
$ws.Range("E3").Value = "  +4.10%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "190.70"
$ws.Range("E5").Value = "  +3.50%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "593.72"
$ws.Range("E6").Value = "  +2.32%  "

# Row 7
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.609"
$ws.Range("E8").Value = "  +0.55%  "

# Row 9
$ws.Range("E9").Value = "  +2.07%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.76"
$ws.Range("E10").Value = "  +2.74%  "

# Row 11
$ws.Range("E11").Value = "  +1.77%  "

# Row 12
$ws.Range("D12").Value = "3.976.01"
$ws.Range("E12").Value = "  +4.64%  "

# Row 13
$ws.Range("E13").Value = "  -0.87%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.70"
$ws.Range("E14").Value = "  +3.70%  "

# Row 15
$ws.Range("D15").Value = "69.609.67"

# Row 16
$ws.Range("E16").Value = "  +1.46%  "

# Row 17
$ws.Range("D17").Value = "3.387.24"
$ws.Range("E17").Value = "  +5.39%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "451.64"
$ws.Range("E18").Value = "  +14.46%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.83"
$ws.Range("E19").Value = "  +1.17%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.81"
$ws.Range("E20").Value = "  +1.95%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.82"
$ws.Range("E21").Value = "  +2.98%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "76.22"
$ws.Range("E22").Value = "  +6.68%  "

# Row 23
$ws.Range("E23").Value = "  -0.05%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.522"
$ws.Range("E24").Value = "  +1.00%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000122"
$ws.Range("E25").Value = "  +3.80%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.192"
$ws.Range("E26").Value = "  +2.64%  "

# Row 27
$ws.Range("E27").Value = "  -0.95%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -0.13%  "

# Row 29
$ws.Range("E29").Value = "  +3.42%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "23.47"
$ws.Range("E30").Value = "  +3.55%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.64"
$ws.Range("E31").Value = "  +1.88%  "

# Row 32
$ws.Range("E32").Value = "  +2.25%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.98"
$ws.Range("E33").Value = "  -0.42%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.58"
$ws.Range("E35").Value = "  +6.90%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.25"
$ws.Range("E36").Value = "  +2.26%  "

# Row 37
$ws.Range("E37").Value = "  +2.56%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "28.32"
$ws.Range("E38").Value = "  +6.15%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.815"
$ws.Range("E39").Value = "  +0.93%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.60"
$ws.Range("E40").Value = "  +1.33%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.60"
$ws.Range("E41").Value = "  +2.01%  "

# Row 42
$ws.Range("D42").Value = "2.754.64"
$ws.Range("E42").Value = "  +5.31%  "

# Row 43
$ws.Range("E43").Value = "  +1.64%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.57"
$ws.Range("E44").Value = "  +3.35%  "

# Row 45
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.16"
$ws.Range("E45").Value = "  +1.17%  "

# Row 46
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0688"
$ws.Range("E46").Value = "  -0.01%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "340.00"
$ws.Range("E47").Value = "  +1.70%  "

# Row 48
$ws.Range("E48").Value = "  +2.50%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.74"
$ws.Range("E49").Value = "  +6.34%  "

# Row 50
$ws.Range("E50").Value = "  +5.15%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.33"
$ws.Range("E51").Value = "  -0.03%  "
